$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both carry the same three updated stats
# (F2: 5779 -> 5798, F5: 974 -> 976, F6: 62 -> 66).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5798
    $ws.Range("F5").Value = 976
    $ws.Range("F6").Value = 66
}
